$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "El servidor recibirá mensajes desde cada uno de los virlocs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El servidor recibirá mensajes desde cada uno de los virlocs.", 2)

$d.Content.Find.Execute(
    "El mensaje contendrá un ID de dispositivo, un tipo de mensaje, un numero de mensaje y un checksum.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El mensaje contendrá un ID de dispositivo, un tipo de mensaje, un numero de mensaje y un checksum.", 2)

$d.Content.Find.Execute(
    "El servidor, una vez recibido el mensaje deberá procesarlo y responder al virloc un mensaje de confirmación.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El servidor, una vez recibido el mensaje deberá procesarlo y responder al virloc un mensaje de confirmación.", 2)
